# "invio (e creazione) pdf ordine alla consegna"
#
# 1) Paragraph "SEO google e cazzi." -> split into several runs (with a
#    spell-check proofErr wrap around "google") and extend the sentence
#    with " (Attendi dominio)." ; the _GoBack bookmark is relocated inside
#    this new text (right after "dominio").
# 2) Delete the "Rivedi caricamento immagini." bullet entirely.
# 3) The paragraph that used to hold the (now relocated) _GoBack bookmark
#    becomes an empty paragraph (handled automatically by moving the
#    bookmark in step 1, since a document may only have one "_GoBack").
# 4) Insert a new empty heading-style paragraph (orange, bold, size 36,
#    underlined) right after the Giacomo/Francesco-list separator blank
#    paragraphs.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: locate & rewrite the "SEO google e cazzi." paragraph
# ---------------------------------------------------------------------
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -eq "SEO google e cazzi.`r") {
        $found = $true

        $pStart = $para.Range.Start
        $pEnd = $para.Range.End
        $textRange = $d.Range($pStart, $pEnd - 1)

        $newXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + `
            '<w:r><w:rPr><w:b/><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve">SEO </w:t></w:r>' + `
            '<w:proofErr w:type="spellStart"/>' + `
            '<w:r><w:rPr><w:b/><w:color w:val="FF0000"/></w:rPr><w:t>google</w:t></w:r>' + `
            '<w:proofErr w:type="spellEnd"/>' + `
            '<w:r><w:rPr><w:b/><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve"> e cazzi</w:t></w:r>' + `
            '<w:r><w:rPr><w:b/><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve"> (Attendi dominio</w:t></w:r>' + `
            '<w:r><w:rPr><w:b/><w:color w:val="FF0000"/></w:rPr><w:t>)</w:t></w:r>' + `
            '<w:r><w:rPr><w:b/><w:color w:val="FF0000"/></w:rPr><w:t>.</w:t></w:r>' + `
            '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $textRange.InsertXML($newXml)

        # Re-fetch the (now longer) paragraph and drop the _GoBack bookmark
        # right after "... (Attendi dominio", i.e. just before the ")".
        $rewritten = $d.Paragraphs($i).Range
        $offset = $rewritten.Text.IndexOf("dominio") + 7
        $bmPos = $rewritten.Start + $offset
        $bmRange = $d.Range($bmPos, $bmPos)
        $d.Bookmarks.Add("_GoBack", $bmRange)

        break
    }
}
if (-not $found) {
    throw "Could not find paragraph 'SEO google e cazzi.'"
}

# ---------------------------------------------------------------------
# Step 2: delete the "Rivedi caricamento immagini." bullet
# ---------------------------------------------------------------------
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -eq "Rivedi caricamento immagini.`r") {
        $para.Range.Delete()
        $found = $true
        break
    }
}
if (-not $found) {
    throw "Could not find paragraph 'Rivedi caricamento immagini.'"
}

# ---------------------------------------------------------------------
# Step 3: insert a new empty heading paragraph (orange/bold/36pt/underline)
#          right before the first of the existing underlined separators.
# ---------------------------------------------------------------------
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    $rng = $para.Range
    if ($rng.Text -eq "`r" -and $rng.Font.Underline -ne 0 -and $rng.Font.Bold -eq -1 -and $rng.Font.Size -eq 18) {
        $insPos = $rng.Start
        $insRange = $d.Range($insPos, $insPos)
        $newParaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
            '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:b/><w:color w:val="F79646" w:themeColor="accent6"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:u w:val="single"/></w:rPr></w:pPr></w:p>' + `
            '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $insRange.InsertXML($newParaXml)
        $found = $true
        break
    }
}
if (-not $found) {
    throw "Could not find the underlined separator paragraph to insert before"
}
